$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells with refreshed doc names / file revisions
$ws.Range("D3").Value = "LONGULT.txt, MEDU.txt, SHORTU.txt"
$ws.Range("D4").Value = "CMS32_DESC_LONG_SHORT_DX.txt"
$ws.Range("D6").Value = "icd10cm_order_2017.txt"
$ws.Range("D7").Value = "icd10pcs_order_2017.txt"
$ws.Range("D10").Value = "sct2_Description_Snapshot-en_US1000124_20170901.txt"

# Add new row for HCPCS (details tbd when loader is written)
$ws.Range("A11").Value = "HCPCS"
$ws.Range("B11").Value = "txt"
$ws.Range("C11").Value = "tab separated"
$ws.Range("D11").Value = "HCPC2018_CONTR_ANWEB.txt"

$ws.Range("E11").Select()
